$wb = $excel.ActiveWorkbook

# --- Sheet "Tabelle1": append a new tracked entry (row 53) ---
$ws1 = $wb.Worksheets.Item("Tabelle1")

$ws1.Range("A53").Value = 45358
$ws1.Range("A53").NumberFormat = "d-mmm"
$ws1.Range("B53").Value = 2
$ws1.Range("C53").Value = "Proposal"

# --- Sheet "Zeitplan": move the "schreiben/Website/Features/Spazi" block
#     from E10:F13 to B7:C10, and add the new timeline notes in A12:A17 ---
$ws2 = $wb.Worksheets.Item("Zeitplan")

$ws2.Range("E10:F13").ClearContents()

$ws2.Range("B7").Value = "schreiben"
$ws2.Range("C7").Value = "2 Monate"
$ws2.Range("B8").Value = "Features jeweils 1 Monat"
$ws2.Range("B9").Value = "Website "
$ws2.Range("C9").Value = "2 Monate"
$ws2.Range("B10").Value = "1 Monat Spazi"

$ws2.Range("A12").Value = "Stand Anfang März"
$ws2.Range("A13").Value = "Website fertig + why shapes were added Ende - März"
$ws2.Range("A14").Value = "für jeden Algorithmus 1 Monat + 1 Monat Spazi = Ende Juli"
$ws2.Range("A15").Value = "Alles auf VM, testen, Experimente 2 Monate - Ende September"
$ws2.Range("A16").Value = "Experteninterviews 1 Monat - Ende Oktober"
$ws2.Range("A17").Value = "Schreiben: 2 Monate Ende 2024"

# --- Update selections to match the saved UI state ---
$null = $ws2.Activate()
$null = $ws2.Range("E16").Select()

$null = $ws1.Activate()
$null = $ws1.Range("C53").Select()
